$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.334.27"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "2.655.84"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +0.51%  "

$ws.Range("D9").Value = "2.654.78"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("E13").Value = "  +4.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.54%  "

$ws.Range("D15").Value = "3.130.29"
$ws.Range("E15").Value = "  +3.83%  "

$ws.Range("D16").Value = "63.204.60"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000145"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "

$ws.Range("D18").Value = "2.656.29"
$ws.Range("E18").Value = "  +3.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("E21").Value = "  +3.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.79%  "

$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.67%  "

$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "545.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +18.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("E32").Value = "  +6.87%  "

$ws.Range("E33").Value = "  +9.18%  "

$ws.Range("D34").Value = "0.0₃0808"
$ws.Range("E34").Value = "  +3.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.36%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "

$ws.Range("E40").Value = "  +11.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.77%  "

$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0576"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.75%  "

$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("E48").Value = "  +3.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.45%  "

$ws.Range("E50").Value = "  +7.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
